$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.056.14"
$ws.Cells.Item(2, 5).Value = "  -0.37%  "

$ws.Cells.Item(3, 4).Value = "1.828.32"
$ws.Cells.Item(3, 5).Value = "  -0.31%  "

$ws.Cells.Item(4, 4).Value = "0.9987"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

$ws.Cells.Item(5, 4).Value = "240.63"
$ws.Cells.Item(5, 5).Value = "  -0.32%  "

$ws.Cells.Item(6, 4).Value = "0.6191"
$ws.Cells.Item(6, 5).Value = "  -6.87%  "

$ws.Cells.Item(7, 4).Value = "1.000"
$ws.Cells.Item(7, 5).Value = "  +0.01%  "

$ws.Cells.Item(8, 4).Value = "44.46"
$ws.Cells.Item(8, 5).Value = "  +6.34%  "

$ws.Cells.Item(9, 4).Value = "0.07345"
$ws.Cells.Item(9, 5).Value = "  -1.03%  "

$ws.Cells.Item(10, 4).Value = "0.2917"
$ws.Cells.Item(10, 5).Value = "  -0.69%  "

$ws.Cells.Item(11, 4).Value = "22.69"
$ws.Cells.Item(11, 5).Value = "  +0.16%  "

$ws.Cells.Item(12, 4).Value = "0.07685"
$ws.Cells.Item(12, 5).Value = "  -0.63%  "

$ws.Cells.Item(13, 4).Value = "1.826.92"
$ws.Cells.Item(13, 5).Value = "  +2.77%  "

$ws.Cells.Item(14, 4).Value = "4.965"
$ws.Cells.Item(14, 5).Value = "  -0.47%  "

$ws.Cells.Item(15, 4).Value = "0.6626"
$ws.Cells.Item(15, 5).Value = "  -1.02%  "

$ws.Cells.Item(16, 4).Value = "81.85"
$ws.Cells.Item(16, 5).Value = "  -1.42%  "

$ws.Cells.Item(17, 4).Value = "0.000009006"
$ws.Cells.Item(17, 5).Value = "  +7.29%  "

$ws.Cells.Item(18, 4).Value = "6.029"
$ws.Cells.Item(18, 5).Value = "  -1.17%  "

$ws.Cells.Item(19, 4).Value = "29.047.09"
$ws.Cells.Item(19, 5).Value = "  -0.26%  "

$ws.Cells.Item(20, 4).Value = "2.074.38"
$ws.Cells.Item(20, 5).Value = "  +0.24%  "

$ws.Cells.Item(21, 4).Value = "225.56"
$ws.Cells.Item(21, 5).Value = "  -0.76%  "

$ws.Cells.Item(22, 4).Value = "12.36"
$ws.Cells.Item(22, 5).Value = "  -0.94%  "

$ws.Cells.Item(23, 4).Value = "1.000"
$ws.Cells.Item(23, 5).Value = "  -0.08%  "

$ws.Cells.Item(24, 4).Value = "7.138"
$ws.Cells.Item(24, 5).Value = "  -0.65%  "

$ws.Cells.Item(25, 4).Value = "1.000"
$ws.Cells.Item(25, 5).Value = "  +0.00%  "

$ws.Cells.Item(26, 4).Value = "160.07"
$ws.Cells.Item(26, 5).Value = "  +0.21%  "

$ws.Cells.Item(27, 4).Value = "8.433"
$ws.Cells.Item(27, 5).Value = "  -2.29%  "

$ws.Cells.Item(28, 4).Value = "0.1354"
$ws.Cells.Item(28, 5).Value = "  -3.85%  "

$ws.Cells.Item(29, 4).Value = "17.79"
$ws.Cells.Item(29, 5).Value = "  -0.90%  "

$ws.Cells.Item(30, 4).Value = "1.493"
$ws.Cells.Item(30, 5).Value = "  -1.18%  "

$ws.Cells.Item(31, 4).Value = "4.042"
$ws.Cells.Item(31, 5).Value = "  -0.12%  "

$ws.Cells.Item(32, 4).Value = "4.053"
$ws.Cells.Item(32, 5).Value = "  -1.56%  "

$ws.Cells.Item(33, 4).Value = "1.199"
$ws.Cells.Item(33, 5).Value = "  +0.40%  "

$ws.Cells.Item(34, 4).Value = "0.05265"
$ws.Cells.Item(34, 5).Value = "  -1.22%  "

$ws.Cells.Item(35, 4).Value = "1.838"
$ws.Cells.Item(35, 5).Value = "  -1.99%  "

$ws.Cells.Item(36, 4).Value = "1.150"
$ws.Cells.Item(36, 5).Value = "  +1.10%  "

$ws.Cells.Item(37, 4).Value = "0.7308"
$ws.Cells.Item(37, 5).Value = "  -3.55%  "

$ws.Cells.Item(38, 4).Value = "2.647"
$ws.Cells.Item(38, 5).Value = "  -0.93%  "

$ws.Cells.Item(39, 4).Value = "1.297.94"
$ws.Cells.Item(39, 5).Value = "  +2.03%  "

$ws.Cells.Item(40, 4).Value = "2.749"
$ws.Cells.Item(40, 5).Value = "  +0.56%  "

$ws.Cells.Item(41, 4).Value = "0.01783"
$ws.Cells.Item(41, 5).Value = "  -0.89%  "

$ws.Cells.Item(42, 4).Value = "6.302"
$ws.Cells.Item(42, 5).Value = "  +5.53%  "

$ws.Cells.Item(43, 4).Value = "0.8991"
$ws.Cells.Item(43, 5).Value = "  -3.28%  "

$ws.Cells.Item(44, 4).Value = "0.9990"
$ws.Cells.Item(44, 5).Value = "  -0.32%  "

$ws.Cells.Item(45, 4).Value = "101.92"
$ws.Cells.Item(45, 5).Value = "  -0.87%  "

$ws.Cells.Item(46, 4).Value = "1.972.40"
$ws.Cells.Item(46, 5).Value = "  +0.35%  "

$ws.Cells.Item(47, 2).Value = "Aave"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(47, 4).Value = "63.96"
$ws.Cells.Item(47, 5).Value = "  +0.99%  "

$ws.Cells.Item(48, 4).Value = "0.5113"
$ws.Cells.Item(48, 5).Value = "  -0.94%  "

$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(49, 4).Value = "0.00000000120"
$ws.Cells.Item(49, 5).Value = "  -0.10%  "

$ws.Cells.Item(50, 4).Value = "1.717"
$ws.Cells.Item(50, 5).Value = "  -3.15%  "

$ws.Cells.Item(51, 4).Value = "0.3965"
$ws.Cells.Item(51, 5).Value = "  -1.78%  "
